$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.968.00"
$ws.Range("D3").Value = "1.655.32"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("D12").Value = "1.889.11"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").Value = "1.650.49"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "26.980.53"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "0.0₃0739"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.26"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.85"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "1.552.74"
$ws.Range("E32").Value = "  +3.67%  "
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("E35").Value = "  +9.22%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("E38").Value = "  +9.04%  "
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.71%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +8.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").Value = "1.796.98"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0993"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.35%  "
